$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ADC Module (row 3) and 70CM Beacon (row 2) status updated to "Complete"
$ws.Range("C2").Value = "Complete"
$ws.Range("C3").Value = "Complete"

# Update the view: scroll back to top, select C4 instead of C22
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C4").Select()
